$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1926803013993541
$ws.Range("C2").Value = 0.5640473627556513
$ws.Range("J2").Value = 0.02368137782561894
$ws.Range("P2").Value = 0.1313240043057051
$ws.Range("S2").Value = 0.08826695371367062
$ws.Range("B3").Value = 0.009107468123861567
$ws.Range("C3").Value = 0.03096539162112933
$ws.Range("J3").Value = 0.0273224043715847
$ws.Range("P3").Value = 0.7522768670309654
$ws.Range("S3").Value = 0.180327868852459
$ws.Range("J4").Value = 0.03424657534246575
$ws.Range("P4").Value = 0.7191780821917808
$ws.Range("S4").Value = 0.2465753424657534
$ws.Range("B6").Value = 0.05854430379746835
$ws.Range("D6").Value = 0.01424050632911392
$ws.Range("F6").Value = 0.07911392405063292
$ws.Range("J6").Value = 0.2357594936708861
$ws.Range("O6").Value = 0.0189873417721519
$ws.Range("Q6").Value = 0.1645569620253164
$ws.Range("R6").Value = 0.06487341772151899
$ws.Range("S6").Value = 0.3639240506329114
$ws.Range("B7").Value = 0.1092985318107667
$ws.Range("D7").Value = 0.02446982055464927
$ws.Range("E7").Value = 0.001631321370309951
$ws.Range("F7").Value = 0.05383360522022838
$ws.Range("J7").Value = 0.1402936378466558
$ws.Range("O7").Value = 0.02283849918433932
$ws.Range("Q7").Value = 0.1908646003262643
$ws.Range("R7").Value = 0.0701468189233279
$ws.Range("S7").Value = 0.3866231647634584
$ws.Range("B8").Value = 0.09672505712109672
$ws.Range("D8").Value = 0.01294744859101295
$ws.Range("F8").Value = 0.0594059405940594
$ws.Range("J8").Value = 0.1111957349581112
$ws.Range("O8").Value = 0.02132520944402132
$ws.Range("Q8").Value = 0.1751713632901752
$ws.Range("R8").Value = 0.1089108910891089
$ws.Range("S8").Value = 0.4143183549124143
$ws.Range("B9").Value = 0.09404388714733543
$ws.Range("D9").Value = 0.02821316614420063
$ws.Range("E9").Value = 0.001567398119122257
$ws.Range("F9").Value = 0.054858934169279
$ws.Range("J9").Value = 0.1442006269592477
$ws.Range("O9").Value = 0.01567398119122257
$ws.Range("Q9").Value = 0.170846394984326
$ws.Range("R9").Value = 0.07680250783699059
$ws.Range("S9").Value = 0.4137931034482759
$ws.Range("B10").Value = 0.109201213346815
$ws.Range("D10").Value = 0.02199191102123357
$ws.Range("E10").Value = 0.001011122345803842
$ws.Range("F10").Value = 0.05915065722952477
$ws.Range("J10").Value = 0.1395348837209302
$ws.Range("O10").Value = 0.0166835187057634
$ws.Range("Q10").Value = 0.2143579373104146
$ws.Range("R10").Value = 0.08190091001011122
$ws.Range("S10").Value = 0.3561678463094035
$ws.Range("G11").Value = 0.1445916114790287
$ws.Range("J11").Value = 0.07947019867549669
$ws.Range("K11").Value = 0.1920529801324503
$ws.Range("L11").Value = 0.5750551876379691
$ws.Range("S11").Value = 0.008830022075055188
$ws.Range("F12").Value = 0.001838235294117647
$ws.Range("G12").Value = 0.7463235294117647
$ws.Range("J12").Value = 0.1746323529411765
$ws.Range("K12").Value = 0.005514705882352942
$ws.Range("L12").Value = 0.04044117647058824
$ws.Range("S12").Value = 0.03125
$ws.Range("F13").Value = 0.007462686567164179
$ws.Range("G13").Value = 0.6716417910447762
$ws.Range("J13").Value = 0.2686567164179104
$ws.Range("S13").Value = 0.05223880597014925
$ws.Range("G14").Value = 0.8
$ws.Range("J14").Value = 0.2
$ws.Range("F15").Value = 0.0256797583081571
$ws.Range("H15").Value = 0.1359516616314199
$ws.Range("I15").Value = 0.06646525679758308
$ws.Range("J15").Value = 0.3776435045317221
$ws.Range("K15").Value = 0.07401812688821752
$ws.Range("M15").Value = 0.01661631419939577
$ws.Range("O15").Value = 0.08006042296072508
$ws.Range("S15").Value = 0.2235649546827795
$ws.Range("F16").Value = 0.01470588235294118
$ws.Range("H16").Value = 0.184640522875817
$ws.Range("I16").Value = 0.08006535947712418
$ws.Range("J16").Value = 0.4313725490196079
$ws.Range("K16").Value = 0.09967320261437909
$ws.Range("M16").Value = 0.01797385620915033
$ws.Range("N16").Value = 0.001633986928104575
$ws.Range("O16").Value = 0.05882352941176471
$ws.Range("S16").Value = 0.1111111111111111
$ws.Range("F17").Value = 0.01717967072297781
$ws.Range("H17").Value = 0.1703650680028633
$ws.Range("I17").Value = 0.101646385110952
$ws.Range("J17").Value = 0.4022906227630637
$ws.Range("K17").Value = 0.09878310665712241
$ws.Range("M17").Value = 0.0164638511095204
$ws.Range("N17").Value = 0.001431639226914817
$ws.Range("O17").Value = 0.06943450250536864
$ws.Range("S17").Value = 0.1224051539012169
$ws.Range("F18").Value = 0.02317880794701987
$ws.Range("H18").Value = 0.2036423841059603
$ws.Range("I18").Value = 0.07450331125827815
$ws.Range("J18").Value = 0.4337748344370861
$ws.Range("K18").Value = 0.09105960264900662
$ws.Range("M18").Value = 0.02317880794701987
$ws.Range("N18").Value = 0.001655629139072848
$ws.Range("O18").Value = 0.04304635761589404
$ws.Range("S18").Value = 0.1059602649006623
$ws.Range("F19").Value = 0.01663146779303062
$ws.Range("H19").Value = 0.1979936642027455
$ws.Range("I19").Value = 0.09371700105596621
$ws.Range("J19").Value = 0.3743400211193242
$ws.Range("K19").Value = 0.1127243928194298
$ws.Range("M19").Value = 0.02006335797254488
$ws.Range("N19").Value = 0.0002639915522703273
$ws.Range("O19").Value = 0.06256599788806758
$ws.Range("S19").Value = 0.1217001055966209
